$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Summer 24 week 9 inputs - update matchup averages
$ws.Range("C2").Value = 1.29
$ws.Range("B3").Value = 1.55
$ws.Range("E3").Value = 1.32
$ws.Range("E4").Value = 1.22
$ws.Range("C5").Value = 1.35
$ws.Range("D5").Value = 1.34
$ws.Range("F5").Value = 1.05
$ws.Range("E7").Value = 1.89
